$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: insert 7 new blank rows at 17:23 (shifts old row 17.. down by 7)
$ws.Rows("17:23").Insert()

# Step 2: insert 2 new blank rows at 25:26 (old row 17 data now sits at row 24;
# this pushes old rows 18-25, now at 25-32, further down to 27-34)
$ws.Rows("25:26").Insert()

# Fill newly inserted rows 17-23
# Row 17
$ws.Range("A17").NumberFormat = "@"
$ws.Range("A17").Value = '803608208'
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = '2/24/2025'
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = 'LA PAMPA 1001'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '803608208'
$ws.Range("F17").NumberFormat = "@"
$ws.Range("F17").Value = 'INCO'
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = 'Pendiente'
$ws.Range("H17").NumberFormat = "@"
$ws.Range("H17").Value = 'Cambiar columna base corroida prioridad media '
$ws.Range("I17").Value = 0
$ws.Range("J17").NumberFormat = "@"
$ws.Range("J17").Value = 'Cambio'
$ws.Range("K17").NumberFormat = "@"
$ws.Range("K17").Value = 'Sin equipos'
$ws.Range("L17").NumberFormat = "@"
$ws.Range("L17").Value = 'Pasante'
$ws.Range("M17").Value = -58.439727
$ws.Range("N17").Value = -34.556261
$ws.Range("O17").NumberFormat = "@"
$ws.Range("O17").Value = 'Saavedra'
$ws.Range("P17").NumberFormat = "@"
$ws.Range("P17").Value = 'Capital Norte'

# Row 18
$ws.Range("A18").NumberFormat = "@"
$ws.Range("A18").Value = '803608178'
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = '2/24/2025'
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = 'HERNANDEZ JOSE 1451'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '803608178'
$ws.Range("F18").NumberFormat = "@"
$ws.Range("F18").Value = 'INCO'
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = 'Pendiente'
$ws.Range("H18").NumberFormat = "@"
$ws.Range("H18").Value = 'Cambiar columna 114 y efectuar transferencias base corroida prioridad media '
$ws.Range("I18").Value = 0
$ws.Range("J18").NumberFormat = "@"
$ws.Range("J18").Value = 'Cambio'
$ws.Range("K18").NumberFormat = "@"
$ws.Range("K18").Value = 'Nodo/Fuente Teco'
$ws.Range("L18").NumberFormat = "@"
$ws.Range("L18").Value = 'Pasante'
$ws.Range("M18").Value = -58.443936
$ws.Range("N18").Value = -34.560145
$ws.Range("O18").NumberFormat = "@"
$ws.Range("O18").Value = 'Saavedra'
$ws.Range("P18").NumberFormat = "@"
$ws.Range("P18").Value = 'Capital Norte'

# Row 19
$ws.Range("A19").NumberFormat = "@"
$ws.Range("A19").Value = '804427439'
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = '4/1/2025'
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = 'Conesa 2195'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '804427439'
$ws.Range("F19").NumberFormat = "@"
$ws.Range("F19").Value = 'INCO'
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = 'Pendiente'
$ws.Range("H19").NumberFormat = "@"
$ws.Range("H19").Value = 'Picada'
$ws.Range("I19").Value = 0
$ws.Range("J19").NumberFormat = "@"
$ws.Range("J19").Value = 'Cambio'
$ws.Range("K19").NumberFormat = "@"
$ws.Range("K19").Value = 'Sin equipos'
$ws.Range("L19").NumberFormat = "@"
$ws.Range("L19").Value = 'Pasante'
$ws.Range("M19").Value = -58.463015
$ws.Range("N19").Value = -34.564505
$ws.Range("O19").NumberFormat = "@"
$ws.Range("O19").Value = 'Colegiales'
$ws.Range("P19").NumberFormat = "@"
$ws.Range("P19").Value = 'Capital Norte'

# Row 20
$ws.Range("A20").NumberFormat = "@"
$ws.Range("A20").Value = '804568979'
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = '4/8/2025'
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = 'Quesada 2710'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '804568979'
$ws.Range("F20").NumberFormat = "@"
$ws.Range("F20").Value = 'INCO'
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = 'Pendiente'
$ws.Range("H20").NumberFormat = "@"
$ws.Range("H20").Value = 'Picada'
$ws.Range("I20").Value = 1
$ws.Range("J20").NumberFormat = "@"
$ws.Range("J20").Value = 'Cambio'
$ws.Range("K20").NumberFormat = "@"
$ws.Range("K20").Value = 'Sin equipos'
$ws.Range("L20").NumberFormat = "@"
$ws.Range("L20").Value = 'Pasante'
$ws.Range("M20").Value = -58.466348
$ws.Range("N20").Value = -34.556028
$ws.Range("O20").NumberFormat = "@"
$ws.Range("O20").Value = 'Saavedra'
$ws.Range("P20").NumberFormat = "@"
$ws.Range("P20").Value = 'Capital Norte'

# Row 21
$ws.Range("A21").NumberFormat = "@"
$ws.Range("A21").Value = '805507192'
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = '4/28/2025'
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = 'Virrey Arredondo 2821'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '805507192'
$ws.Range("F21").NumberFormat = "@"
$ws.Range("F21").Value = 'INCO'
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = 'Pendiente'
$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value = 'Picada'
$ws.Range("I21").Value = 1
$ws.Range("J21").NumberFormat = "@"
$ws.Range("J21").Value = 'Cambio'
$ws.Range("K21").NumberFormat = "@"
$ws.Range("K21").Value = 'Sin equipos'
$ws.Range("L21").NumberFormat = "@"
$ws.Range("L21").Value = 'Terminal'
$ws.Range("M21").Value = -58.454065
$ws.Range("N21").Value = -34.57105
$ws.Range("O21").NumberFormat = "@"
$ws.Range("O21").Value = 'Colegiales'
$ws.Range("P21").NumberFormat = "@"
$ws.Range("P21").Value = 'Capital Norte'

# Row 22
$ws.Range("A22").NumberFormat = "@"
$ws.Range("A22").Value = '805655355'
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = '5/5/2025'
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = 'Arce 867'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '14'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '805655355'
$ws.Range("F22").NumberFormat = "@"
$ws.Range("F22").Value = 'INCO'
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = 'Pendiente'
$ws.Range("H22").NumberFormat = "@"
$ws.Range("H22").Value = 'Picada'
$ws.Range("I22").Value = 1
$ws.Range("J22").NumberFormat = "@"
$ws.Range("J22").Value = 'Cambio'
$ws.Range("K22").NumberFormat = "@"
$ws.Range("K22").Value = 'Sin equipos'
$ws.Range("L22").NumberFormat = "@"
$ws.Range("L22").Value = 'Pasante'
$ws.Range("M22").Value = -58.436255
$ws.Range("N22").Value = -34.567733
$ws.Range("O22").NumberFormat = "@"
$ws.Range("O22").Value = 'Palermo'
$ws.Range("P22").NumberFormat = "@"
$ws.Range("P22").Value = 'Capital Sur'

# Row 23
$ws.Range("A23").NumberFormat = "@"
$ws.Range("A23").Value = '805655369'
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = '5/5/2025'
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = 'Benjamin Matienzo 1524'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '805655369'
$ws.Range("F23").NumberFormat = "@"
$ws.Range("F23").Value = 'INCO'
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = 'Pendiente'
$ws.Range("H23").NumberFormat = "@"
$ws.Range("H23").Value = 'Picada'
$ws.Range("I23").Value = 1
$ws.Range("J23").NumberFormat = "@"
$ws.Range("J23").Value = 'Cambio'
$ws.Range("K23").NumberFormat = "@"
$ws.Range("K23").Value = 'Sin equipos'
$ws.Range("L23").NumberFormat = "@"
$ws.Range("L23").Value = 'Terminal'
$ws.Range("M23").Value = -58.43247
$ws.Range("N23").Value = -34.566492
$ws.Range("O23").NumberFormat = "@"
$ws.Range("O23").Value = 'Palermo'
$ws.Range("P23").NumberFormat = "@"
$ws.Range("P23").Value = 'Capital Sur'

# Fill newly inserted rows 25-26
# Row 25
$ws.Range("A25").NumberFormat = "@"
$ws.Range("A25").Value = '805707245'
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = '5/6/2025'
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value = 'Soldado de la Independencia 1298'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '14'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '805707245'
$ws.Range("F25").NumberFormat = "@"
$ws.Range("F25").Value = 'INCO'
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = 'Pendiente'
$ws.Range("H25").NumberFormat = "@"
$ws.Range("H25").Value = 'Picada - Con fuente teco'
$ws.Range("I25").Value = 1
$ws.Range("J25").NumberFormat = "@"
$ws.Range("J25").Value = 'Cambio'
$ws.Range("K25").NumberFormat = "@"
$ws.Range("K25").Value = 'Fuente Teco'
$ws.Range("L25").NumberFormat = "@"
$ws.Range("L25").Value = 'Pasante'
$ws.Range("M25").Value = -58.440507
$ws.Range("N25").Value = -34.564016
$ws.Range("O25").NumberFormat = "@"
$ws.Range("O25").Value = 'Colegiales'
$ws.Range("P25").NumberFormat = "@"
$ws.Range("P25").Value = 'Capital Norte'

# Row 26
$ws.Range("A26").NumberFormat = "@"
$ws.Range("A26").Value = '805722772'
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = '5/7/2025'
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = 'Luis Maria Campos 1336'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '14'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '805722772'
$ws.Range("F26").NumberFormat = "@"
$ws.Range("F26").Value = 'INCO'
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = 'Pendiente'
$ws.Range("H26").NumberFormat = "@"
$ws.Range("H26").Value = 'Picada'
$ws.Range("I26").Value = 1
$ws.Range("J26").NumberFormat = "@"
$ws.Range("J26").Value = 'Cambio'
$ws.Range("K26").NumberFormat = "@"
$ws.Range("K26").Value = 'Sin equipos'
$ws.Range("L26").NumberFormat = "@"
$ws.Range("L26").Value = 'Pasante'
$ws.Range("M26").Value = -58.44191
$ws.Range("N26").Value = -34.564245
$ws.Range("O26").NumberFormat = "@"
$ws.Range("O26").Value = 'Colegiales'
$ws.Range("P26").NumberFormat = "@"
$ws.Range("P26").Value = 'Capital Norte'

# Step 3: append 2 new rows at the end (35-36)
# Row 35
$ws.Range("A35").NumberFormat = "@"
$ws.Range("A35").Value = '-523'
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = '7/20/2025'
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = 'Luis Maria Campos 585'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '14'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '808460898'
$ws.Range("F35").NumberFormat = "@"
$ws.Range("F35").Value = 'INCO'
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = 'Pendiente'
$ws.Range("H35").NumberFormat = "@"
$ws.Range("H35").Value = 'Picada'
$ws.Range("I35").Value = 1
$ws.Range("J35").NumberFormat = "@"
$ws.Range("J35").Value = 'Cambio'
$ws.Range("K35").NumberFormat = "@"
$ws.Range("K35").Value = 'Sin equipos'
$ws.Range("L35").NumberFormat = "@"
$ws.Range("L35").Value = 'Pasante'
$ws.Range("M35").Value = -58.434668
$ws.Range("N35").Value = -34.571258
$ws.Range("O35").NumberFormat = "@"
$ws.Range("O35").Value = 'Palermo'
$ws.Range("P35").NumberFormat = "@"
$ws.Range("P35").Value = 'Capital Sur'

# Row 36
$ws.Range("A36").NumberFormat = "@"
$ws.Range("A36").Value = '-524'
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = '7/21/2025'
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = 'Luis Maria Campos 509'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '14'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '808460897'
$ws.Range("F36").NumberFormat = "@"
$ws.Range("F36").Value = 'INCO'
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = 'Pendiente'
$ws.Range("H36").NumberFormat = "@"
$ws.Range("H36").Value = 'Picada'
$ws.Range("I36").Value = 1
$ws.Range("J36").NumberFormat = "@"
$ws.Range("J36").Value = 'Cambio'
$ws.Range("K36").NumberFormat = "@"
$ws.Range("K36").Value = 'Sin equipos'
$ws.Range("L36").NumberFormat = "@"
$ws.Range("L36").Value = 'Pasante'
$ws.Range("M36").Value = -58.434194
$ws.Range("N36").Value = -34.571754
$ws.Range("O36").NumberFormat = "@"
$ws.Range("O36").Value = 'Palermo'
$ws.Range("P36").NumberFormat = "@"
$ws.Range("P36").Value = 'Capital Sur'

